$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.607.05"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "1.584.82"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "1.807.12"
$ws.Range("E12").Value = "  -2.81%  "
$ws.Range("D13").Value = "1.585.60"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "26.599.21"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").Value = "  -3.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.674"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +24.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("D35").Value = "1.326.62"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0172"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.782"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.720.11"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.832"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0986"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
